# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: drop the " stds" suffix (format as text so "2021".."2024"
# stay strings instead of being auto-coerced to numbers by Excel)
$ws.Range("B1:E1").NumberFormat = "@"
$ws.Range("B1").Value = "2021"
$ws.Range("C1").Value = "2022"
$ws.Range("D1").Value = "2023"
$ws.Range("E1").Value = "2024"

# Row 2 (TB) - regenerated std/mean for 2022/2023/2024 columns
$ws.Range("C2").Value = 2.09277672131318
$ws.Range("D2").Value = 2.236972821555524
$ws.Range("E2").Value = 2.10352542763909

# Row 3 (PC)
$ws.Range("C3").Value = 12.97923549891196
$ws.Range("D3").Value = 13.27028011905468
$ws.Range("E3").Value = 12.3708529057527

# Row 4 (dS0)
$ws.Range("C4").Value = 4.056899336615523
$ws.Range("D4").Value = 4.189146948016854
$ws.Range("E4").Value = 3.954634125721075

# Row 5 (dSF)
$ws.Range("C5").Value = 4.260838190045037
$ws.Range("D5").Value = 4.395116360286813
$ws.Range("E5").Value = 4.168971742130561

# Row 6 (K) - now using K instead of Strike#, so all four years are
# recalculated (previously only 2021..2024 "Strike#" values)
$ws.Range("B6").Value = 1.398315329090978
$ws.Range("C6").Value = 1.226182245272708
$ws.Range("D6").Value = 1.265926837711646
$ws.Range("E6").Value = 1.174725659777613

# Row 7 (IP)
$ws.Range("C7").Value = 0.9678892714594651
$ws.Range("D7").Value = 0.9443521491082709
$ws.Range("E7").Value = 0.8979769202998196

# Row 8 (I0)
$ws.Range("C8").Value = 0.967889271459465
$ws.Range("D8").Value = 0.9443521491082709
$ws.Range("E8").Value = 0.8979769202998196
